# Updates cryptos list values (prices & 1h volume %) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.699.05'
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").Value = '2.310.93'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.87'
$ws.Range("E5").Value = '  -0.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.36'
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.503'
$ws.Range("E7").Value = '  -0.42%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.492'
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.16'
$ws.Range("E10").Value = '  -2.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '18.82'
$ws.Range("E11").Value = '  +1.16%  '
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("E14").Value = '  -1.77%  '
$ws.Range("D15").Value = '2.670.25'
$ws.Range("E15").Value = '  +0.65%  '
$ws.Range("D16").Value = '2.318.75'
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.787'
$ws.Range("E17").Value = '  +1.39%  '
$ws.Range("D18").Value = '42.651.85'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.14'
$ws.Range("E19").Value = '  -4.70%  '
$ws.Range("E20").Value = '  +1.86%  '
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("E22").Value = '  +0.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.27'
$ws.Range("E23").Value = '  +5.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.07'
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.41'
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.28'
$ws.Range("E27").Value = '  -1.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.35'
$ws.Range("E28").Value = '  +14.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.07'
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("E30").Value = '  +0.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.05'
$ws.Range("E31").Value = '  -2.42%  '
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("E33").Value = '  +0.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.66'
$ws.Range("E34").Value = '  -1.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.45'
$ws.Range("E35").Value = '  +0.53%  '
$ws.Range("E36").Value = '  +1.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.34'
$ws.Range("E37").Value = '  -0.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.79'
$ws.Range("E38").Value = '  +2.77%  '
$ws.Range("E39").Value = '  -0.30%  '
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.70'
$ws.Range("E40").Value = '  +0.20%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.108'
$ws.Range("E41").Value = '  -0.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.45'
$ws.Range("E42").Value = '  +18.69%  '
$ws.Range("D43").Value = '1.923.10'
$ws.Range("E43").Value = '  -3.72%  '
$ws.Range("E44").Value = '  -0.61%  '
$ws.Range("E45").Value = '  -2.09%  '
$ws.Range("E46").Value = '  -1.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.74'
$ws.Range("E47").Value = '  -0.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.88'
$ws.Range("E48").Value = '  +1.90%  '
$ws.Range("D49").Value = '2.538.79'
$ws.Range("E49").Value = '  +0.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.29'
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.99'
